$d = $word.ActiveDocument

# Locate the heading paragraph that introduces the "serverless" architecture
# section (it contains the word "serverless"). Every paragraph AFTER that
# heading, through the end of the document, needs bold added on top of the
# italic formatting it already carries (both on the paragraph mark and on
# every run), matching the target diff.
$paragraphs = $d.Paragraphs
$count = $paragraphs.Count

$headingIndex = -1
for ($i = 1; $i -le $count; $i++) {
    $para = $paragraphs.Item($i)
    if ($para.Range.Text -like "*serverless*") {
        $headingIndex = $i
    }
}

if ($headingIndex -eq -1) {
    throw "Could not locate the 'serverless' heading paragraph"
}

for ($i = $headingIndex + 1; $i -le $count; $i++) {
    $para = $paragraphs.Item($i)
    $rng = $para.Range
    $rng.Font.Bold = $true
    $rng.Font.BoldBi = $true
}

Write-Output ("Bolded paragraphs " + ($headingIndex + 1) + " through " + $count)
